$d = $word.ActiveDocument
$firstIndex = $d.Paragraphs.Count + 1

# Pre-create all new (empty) paragraphs with their styles first so that
# later italic runs in one paragraph cannot "leak" into the start of the
# next paragraph via formatting inheritance at the insertion point.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Heading1"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Heading2"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Style = "Normal"

# Now fill in the text/formatting for each of those paragraphs.
# --- new paragraph 1 ---
$para = $d.Paragraphs.Item($firstIndex + 0)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Knärot – ekologi samt krav på livsmiljön")

# --- new paragraph 2 ---
$para = $d.Paragraphs.Item($firstIndex + 1)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).")

# --- new paragraph 3 ---
$para = $d.Paragraphs.Item($firstIndex + 2)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Samuel Johnsons doktorsavhandling ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Vidare ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$ins.Font.Italic = $true

# --- new paragraph 4 ---
$para = $d.Paragraphs.Item($firstIndex + 3)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$ins.Font.Italic = $true

# --- new paragraph 5 ---
$para = $d.Paragraphs.Item($firstIndex + 4)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).")

# --- new paragraph 6 ---
$para = $d.Paragraphs.Item($firstIndex + 5)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).")

# --- new paragraph 7 ---
$para = $d.Paragraphs.Item($firstIndex + 6)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Referenser - knärot")

# --- new paragraph 8 ---
$para = $d.Paragraphs.Item($firstIndex + 7)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("de Graaf M & Roberts M.R., 2009. ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Short-term response of the herbaceous layer within leave patches after harvest. ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Forest Ecology and Management 257, 1014-1025")

# --- new paragraph 9 ---
$para = $d.Paragraphs.Item($firstIndex + 8)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Ecological Applications, 22, 2049-2064 ")

# --- new paragraph 10 ---
$para = $d.Paragraphs.Item($firstIndex + 9)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Interactive effects of drought and edge exposure on old-growth forest understory species. ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Landscape Ecology, 37, sid 1839-1853")

# --- new paragraph 11 ---
$para = $d.Paragraphs.Item($firstIndex + 10)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Biological legacies buffer local species extinction after logging. ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Journal of Applied Ecology. 51, 53-62.")

# --- new paragraph 12 ---
$para = $d.Paragraphs.Item($firstIndex + 11)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Skogsstyrelsen, 2022. ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Vägledning för hänsyn till knärot. ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/")

# --- new paragraph 13 ---
$para = $d.Paragraphs.Item($firstIndex + 12)
$pos = $para.Range.Start
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("SLU Artdatabanken, 2021. ")
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Artfaktablad. Naturvård – artfakta. ")
$ins.Font.Italic = $true
$pos = $para.Range.End - 1
$ins = $d.Range($pos, $pos)
$ins.InsertAfter("SLU Artdatabanken, Uppsala ")

# --- update the date stamp in the header ---
foreach ($story in $d.StoryRanges) {
    $story.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
